$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.846.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.862.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "470.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.46%  "
$ws.Range("E7").Value = "  -2.10%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.712"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.163"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000342"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "41.85"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.476.03"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("E14").Value = "  -2.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.921.38"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.83%  "
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.71%  "
$ws.Range("E19").Value = "  -4.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.092.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "426.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("B22").Value = "ImmutableX"
$ws.Range("C22").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.88%  "
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "37.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "722.10"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.36%  "
$ws.Range("E31").Value = "  -5.67%  "
$ws.Range("E32").Value = "  +1.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "41.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0863"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +25.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.05"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.151"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.10%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.28"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.37%  "
$ws.Range("E39").Value = "  -3.13%  "
$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.02"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.37%  "
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.96"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +10.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.340"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.138"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.35%  "
$ws.Range("E46").Value = "  -2.25%  "
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "144.43"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.66%  "
$ws.Range("E51").Value = "  -2.65%  "
